$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for the two new columns ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data values for I2:J37 ---
$I_values = @(9,8,7,8,6,8,7,8,8,9,7,7,5,8,7,7,6,6,6,7,7,6,9,6,8,7,6,6,6,4,8,9,3,7,8,5)
$J_values = @(9,8,7,8,7,8,7,8,8,9,7,7,6,8,7,7,6,7,7,7,7,7,9,7,8,7,6,7,6,5,8,9,4,7,8,5)

for ($idx = 0; $idx -lt $I_values.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $I_values[$idx]
    $ws.Cells.Item($row, 10).Value = $J_values[$idx]
}
